# feat: add 2022-Q1 data
#
# 1. Insert a new sheet "2022-Q1" positioned after "2021-Q4" and before "总计",
#    by cloning the "2021-Q4" sheet (so the fund-holding header row/styles/column
#    layout come along "for free") and then overwriting the data with the
#    2022-Q1 figures (6 funds instead of 3, so 3 extra rows are cloned in too).
# 2. Update the "总计" (totals) sheet: insert a new row on top with the
#    2022-Q1 summary figures, pushing the existing 2021-Q4 / 2021-Q3 rows down.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: build the "2022-Q1" sheet from a copy of "2021-Q4"
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item("2021-Q4")
$template.Copy($null, $template)

$q1 = $wb.Worksheets.Item(3)
$q1.Name = "2022-Q1"

# the template only carries 3 data rows (rows 2-4) - clone row 4's formatting
# down to rows 5-7 so we have the 6 rows of data this quarter needs
$q1.Range("A4:H4").Copy($q1.Range("A5:H5"))
$q1.Range("A4:H4").Copy($q1.Range("A6:H6"))
$q1.Range("A4:H4").Copy($q1.Range("A7:H7"))

# row 2
$q1.Range("A2").Value = 0
$q1.Range("B2").Value = "'001637"
$q1.Range("C2").Value = "嘉实腾讯自选股大数据策略股票"
$q1.Range("D2").Value = "'7.44"
$q1.Range("E2").Value = "'90.10"
$q1.Range("F2").Value = "'1.72"
$q1.Range("G2").Value = "'0.1280"
$q1.Range("H2").Value = 5

# row 3
$q1.Range("A3").Value = 1
$q1.Range("B3").Value = "'011231"
$q1.Range("C3").Value = "光大保德信锦弘混合A"
$q1.Range("D3").Value = "'4.13"
$q1.Range("E3").Value = "'20.96"
$q1.Range("F3").Value = "'0.72"
$q1.Range("G3").Value = "'0.0297"
$q1.Range("H3").Value = 7

# row 4
$q1.Range("A4").Value = 2
$q1.Range("B4").Value = "'011232"
$q1.Range("C4").Value = "光大保德信锦弘混合C"
$q1.Range("D4").Value = "'1.29"
$q1.Range("E4").Value = "'20.96"
$q1.Range("F4").Value = "'0.72"
$q1.Range("G4").Value = "'0.0093"
$q1.Range("H4").Value = 7

# row 5
$q1.Range("A5").Value = 3
$q1.Range("B5").Value = "'005167"
$q1.Range("C5").Value = "嘉实润泽量化一年定期开放混合"
$q1.Range("D5").Value = "'0.56"
$q1.Range("E5").Value = "'27.26"
$q1.Range("F5").Value = "'0.74"
$q1.Range("G5").Value = "'0.0041"
$q1.Range("H5").Value = 1

# row 6
$q1.Range("A6").Value = 4
$q1.Range("B6").Value = "'009918"
$q1.Range("C6").Value = "上银核心成长混合A"
$q1.Range("D6").Value = "'0.13"
$q1.Range("E6").Value = "'91.71"
$q1.Range("F6").Value = "'0.79"
$q1.Range("G6").Value = "'0.0010"
$q1.Range("H6").Value = 8

# row 7
$q1.Range("A7").Value = 5
$q1.Range("B7").Value = "'009919"
$q1.Range("C7").Value = "上银核心成长混合C"
$q1.Range("D7").Value = "'0.07"
$q1.Range("E7").Value = "'91.71"
$q1.Range("F7").Value = "'0.79"
$q1.Range("G7").Value = "'0.0006"
$q1.Range("H7").Value = 8

# ---------------------------------------------------------------------------
# Step 2: update the "总计" sheet with a new 2022-Q1 row on top
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")

$total.Rows(2).Insert()

# restore the index-column style on the freshly inserted row (Insert leaves
# it unstyled) by cloning it from the row right below, then set its value
$total.Range("A3").Copy($total.Range("A2"))
$total.Range("B2:D2").ClearFormats()

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 6
$total.Range("D2").Value = 0.17

# renumber the index column for the rows that shifted down
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
